$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 255.5
$ws.Range("I33").Value = 246.66667
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 246.66667
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -17.66667000000001
$ws.Range("N33").Value = -1058

$ws.Range("H76").Value = 58054.668
$ws.Range("I76").Value = 73871.86
$ws.Range("J76").Value = 2694.5
$ws.Range("K76").Value = 73871.86
$ws.Range("L76").Value = 2694.5
$ws.Range("M76").Value = -73556.86
$ws.Range("N76").Value = -3324.5

$ws.Range("H79").Value = 58054.668
$ws.Range("I79").Value = 73871.86
$ws.Range("J79").Value = 2694.5
$ws.Range("K79").Value = 73871.86
$ws.Range("L79").Value = 2694.5
$ws.Range("M79").Value = -72779.86
$ws.Range("N79").Value = -4878.5

$ws.Range("H129").Value = 18412.457
$ws.Range("I129").Value = 590.5714
$ws.Range("J129").Value = 24214.93
$ws.Range("K129").Value = 1771.7142
$ws.Range("L129").Value = 72644.79000000001
$ws.Range("M129").Value = 3228.2858
$ws.Range("N129").Value = -82644.79000000001

$ws.Range("H132").Value = 3248416.8
$ws.Range("I132").Value = 3403040.2
$ws.Range("K132").Value = 10209120.6
$ws.Range("M132").Value = -10206590.6

$ws.Range("H135").Value = 1490.6364
$ws.Range("I135").Value = 1439.7
$ws.Range("K135").Value = 12957.3
$ws.Range("M135").Value = -10422.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22443.527
$ws.Range("I32").Value = 22850.863
$ws.Range("K32").Value = 22850.863
$ws.Range("M32").Value = -22563.863

$ws.Range("H97").Value = 508.80646
$ws.Range("I97").Value = 477.25
$ws.Range("J97").Value = 803.3333
$ws.Range("K97").Value = 477.25
$ws.Range("L97").Value = 803.3333
$ws.Range("M97").Value = 18.75
$ws.Range("N97").Value = -1795.3333

$ws.Range("H122").Value = 1047.5
$ws.Range("I122").Value = 830.2083
$ws.Range("K122").Value = 2490.6249
$ws.Range("M122").Value = -40.6248999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2721.2
$ws.Range("I86").Value = 1651.5
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 1651.5
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -528.5
$ws.Range("N86").Value = -9246

$ws.Range("H89").Value = 2721.2
$ws.Range("I89").Value = 1651.5
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 8257.5
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -2641.5
$ws.Range("N89").Value = -46232

$ws.Range("H94").Value = 998.36
$ws.Range("I94").Value = 720.5625
$ws.Range("J94").Value = 1492.2222
$ws.Range("K94").Value = 720.5625
$ws.Range("L94").Value = 1492.2222
$ws.Range("M94").Value = -269.5625
$ws.Range("N94").Value = -2394.2222

$ws.Range("H99").Value = 2320.2
$ws.Range("I99").Value = 1926.3334
$ws.Range("J99").Value = 2911
$ws.Range("K99").Value = 1926.3334
$ws.Range("L99").Value = 2911
$ws.Range("M99").Value = -428.3334
$ws.Range("N99").Value = -5907

$ws.Range("H105").Value = 2858
$ws.Range("I105").Value = 2502.875
$ws.Range("J105").Value = 3568.25
$ws.Range("K105").Value = 2502.875
$ws.Range("L105").Value = 3568.25
$ws.Range("M105").Value = -755.875
$ws.Range("N105").Value = -7062.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 379.85715
$ws.Range("I22").Value = 358.72726
$ws.Range("J22").Value = 457.33334
$ws.Range("K22").Value = 358.72726
$ws.Range("L22").Value = 457.33334
$ws.Range("M22").Value = -8.727260000000001
$ws.Range("N22").Value = -1157.33334

$ws.Range("H62").Value = 33335384
$ws.Range("I62").Value = 1745
$ws.Range("J62").Value = 55557810
$ws.Range("K62").Value = 1745
$ws.Range("L62").Value = 55557810
$ws.Range("M62").Value = -1121
$ws.Range("N62").Value = -55559058

$ws.Range("H65").Value = 33335384
$ws.Range("I65").Value = 1745
$ws.Range("J65").Value = 55557810
$ws.Range("K65").Value = 8725
$ws.Range("L65").Value = 277789050
$ws.Range("M65").Value = -5605
$ws.Range("N65").Value = -277795290

$ws.Range("H107").Value = 702.3611
$ws.Range("I107").Value = 809.5789
$ws.Range("J107").Value = 582.5294
$ws.Range("K107").Value = 809.5789
$ws.Range("L107").Value = 582.5294
$ws.Range("M107").Value = 1110.4211
$ws.Range("N107").Value = -4422.5294

$ws.Range("H132").Value = 2634
$ws.Range("I132").Value = 2516.2856
$ws.Range("J132").Value = 2881.2
$ws.Range("K132").Value = 7548.8568
$ws.Range("L132").Value = 8643.599999999999
$ws.Range("M132").Value = -5018.8568
$ws.Range("N132").Value = -13703.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2800
$ws.Range("J49").Value = 2800
$ws.Range("L49").Value = 8400
$ws.Range("N49").Value = -8712

$ws.Range("H131").Value = 1948.5393
$ws.Range("I131").Value = 14835.714
$ws.Range("J131").Value = 848.4146
$ws.Range("K131").Value = 44507.142
$ws.Range("L131").Value = 2545.2438
$ws.Range("M131").Value = -39467.142
$ws.Range("N131").Value = -12625.2438

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 28900634
$ws.Range("I113").Value = 52020344
$ws.Range("K113").Value = 52020344
$ws.Range("M113").Value = -52018174

$ws.Range("H122").Value = 30305822
$ws.Range("I122").Value = 50003070
$ws.Range("J122").Value = 2360
$ws.Range("K122").Value = 150009210
$ws.Range("L122").Value = 7080
$ws.Range("M122").Value = -150006760
$ws.Range("N122").Value = -11980

$ws.Range("H132").Value = 20712.434
$ws.Range("I132").Value = 28667.459
$ws.Range("J132").Value = 2316.4375
$ws.Range("K132").Value = 86002.37699999999
$ws.Range("L132").Value = 6949.3125
$ws.Range("M132").Value = -83472.37699999999
$ws.Range("N132").Value = -12009.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41668908
$ws.Range("I7").Value = 1896.6666
$ws.Range("J7").Value = 55557910
$ws.Range("K7").Value = 1896.6666
$ws.Range("L7").Value = 55557910
$ws.Range("M7").Value = -1784.6666
$ws.Range("N7").Value = -55558134

$ws.Range("H126").Value = 41668908
$ws.Range("I126").Value = 1896.6666
$ws.Range("J126").Value = 55557910
$ws.Range("K126").Value = 5689.9998
$ws.Range("L126").Value = 166673730
$ws.Range("M126").Value = -3219.9998
$ws.Range("N126").Value = -166678670

$ws.Range("H132").Value = 6870.2573
$ws.Range("I132").Value = 10503.579
$ws.Range("K132").Value = 31510.737
$ws.Range("M132").Value = -28980.737

$ws.Range("H133").Value = 24256.5
$ws.Range("J133").Value = 24256.5
$ws.Range("L133").Value = 24256.5
$ws.Range("N133").Value = -29316.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 522.44446
$ws.Range("I100").Value = 562.75
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 1125.5
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = -584.5
$ws.Range("N100").Value = -1482

$ws.Range("H132").Value = 1466.8718
$ws.Range("I132").Value = 755.6667
$ws.Range("J132").Value = 2604.8
$ws.Range("K132").Value = 2267.0001
$ws.Range("L132").Value = 7814.400000000001
$ws.Range("M132").Value = 262.9998999999998
$ws.Range("N132").Value = -12874.4
